# Pievepelago.xlsx update "aggiornamento fino a 8/12" (update through 2021-12-08)
# Appends daily COVID tracking rows 386-464 (dates 44460-44538) to Sheet1,
# extending the sheet dimension from A1:D385 to A1:D464.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the date-column formatting (style index 2: bordered, centered,
# custom date/time number format) from the last existing row (A385) down
# across the whole new block, so the new date cells reuse the existing style
# instead of creating a duplicate one.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new daily data: column A = date serial, B = new positive cases,
# C = 7-day rolling sum, D = 7-day rolling sum per 100k inhabitants.
$ws.Range("A386").Value = 44460
$ws.Range("B386").Value = 0
$ws.Range("C386").Value = 0
$ws.Range("D386").Value = 0
$ws.Range("A387").Value = 44461
$ws.Range("B387").Value = 0
$ws.Range("C387").Value = 0
$ws.Range("D387").Value = 0
$ws.Range("A388").Value = 44462
$ws.Range("B388").Value = 0
$ws.Range("C388").Value = 0
$ws.Range("D388").Value = 0
$ws.Range("A389").Value = 44463
$ws.Range("B389").Value = 2
$ws.Range("C389").Value = 2
$ws.Range("D389").Value = 87.56567425569177
$ws.Range("A390").Value = 44464
$ws.Range("B390").Value = 0
$ws.Range("C390").Value = 2
$ws.Range("D390").Value = 87.56567425569177
$ws.Range("A391").Value = 44465
$ws.Range("B391").Value = 0
$ws.Range("C391").Value = 2
$ws.Range("D391").Value = 87.56567425569177
$ws.Range("A392").Value = 44466
$ws.Range("B392").Value = 0
$ws.Range("C392").Value = 2
$ws.Range("D392").Value = 87.56567425569177
$ws.Range("A393").Value = 44467
$ws.Range("B393").Value = 0
$ws.Range("C393").Value = 2
$ws.Range("D393").Value = 87.56567425569177
$ws.Range("A394").Value = 44468
$ws.Range("B394").Value = 0
$ws.Range("C394").Value = 2
$ws.Range("D394").Value = 87.56567425569177
$ws.Range("A395").Value = 44469
$ws.Range("B395").Value = 0
$ws.Range("C395").Value = 2
$ws.Range("D395").Value = 87.56567425569177
$ws.Range("A396").Value = 44470
$ws.Range("B396").Value = 0
$ws.Range("C396").Value = 0
$ws.Range("D396").Value = 0
$ws.Range("A397").Value = 44471
$ws.Range("B397").Value = 0
$ws.Range("C397").Value = 0
$ws.Range("D397").Value = 0
$ws.Range("A398").Value = 44472
$ws.Range("B398").Value = 0
$ws.Range("C398").Value = 0
$ws.Range("D398").Value = 0
$ws.Range("A399").Value = 44473
$ws.Range("B399").Value = 1
$ws.Range("C399").Value = 1
$ws.Range("D399").Value = 43.78283712784589
$ws.Range("A400").Value = 44474
$ws.Range("B400").Value = 0
$ws.Range("C400").Value = 1
$ws.Range("D400").Value = 43.78283712784589
$ws.Range("A401").Value = 44475
$ws.Range("B401").Value = 0
$ws.Range("C401").Value = 1
$ws.Range("D401").Value = 43.78283712784589
$ws.Range("A402").Value = 44476
$ws.Range("B402").Value = 0
$ws.Range("C402").Value = 1
$ws.Range("D402").Value = 43.78283712784589
$ws.Range("A403").Value = 44477
$ws.Range("B403").Value = 0
$ws.Range("C403").Value = 1
$ws.Range("D403").Value = 43.78283712784589
$ws.Range("A404").Value = 44478
$ws.Range("B404").Value = 0
$ws.Range("C404").Value = 1
$ws.Range("D404").Value = 43.78283712784589
$ws.Range("A405").Value = 44479
$ws.Range("B405").Value = 0
$ws.Range("C405").Value = 1
$ws.Range("D405").Value = 43.78283712784589
$ws.Range("A406").Value = 44480
$ws.Range("B406").Value = 0
$ws.Range("C406").Value = 0
$ws.Range("D406").Value = 0
$ws.Range("A407").Value = 44481
$ws.Range("B407").Value = 0
$ws.Range("C407").Value = 0
$ws.Range("D407").Value = 0
$ws.Range("A408").Value = 44482
$ws.Range("B408").Value = 0
$ws.Range("C408").Value = 0
$ws.Range("D408").Value = 0
$ws.Range("A409").Value = 44483
$ws.Range("B409").Value = 0
$ws.Range("C409").Value = 0
$ws.Range("D409").Value = 0
$ws.Range("A410").Value = 44484
$ws.Range("B410").Value = 0
$ws.Range("C410").Value = 0
$ws.Range("D410").Value = 0
$ws.Range("A411").Value = 44485
$ws.Range("B411").Value = 0
$ws.Range("C411").Value = 0
$ws.Range("D411").Value = 0
$ws.Range("A412").Value = 44486
$ws.Range("B412").Value = 0
$ws.Range("C412").Value = 0
$ws.Range("D412").Value = 0
$ws.Range("A413").Value = 44487
$ws.Range("B413").Value = 0
$ws.Range("C413").Value = 0
$ws.Range("D413").Value = 0
$ws.Range("A414").Value = 44488
$ws.Range("B414").Value = 0
$ws.Range("C414").Value = 0
$ws.Range("D414").Value = 0
$ws.Range("A415").Value = 44489
$ws.Range("B415").Value = 0
$ws.Range("C415").Value = 0
$ws.Range("D415").Value = 0
$ws.Range("A416").Value = 44490
$ws.Range("B416").Value = 0
$ws.Range("C416").Value = 0
$ws.Range("D416").Value = 0
$ws.Range("A417").Value = 44491
$ws.Range("B417").Value = 0
$ws.Range("C417").Value = 0
$ws.Range("D417").Value = 0
$ws.Range("A418").Value = 44492
$ws.Range("B418").Value = 0
$ws.Range("C418").Value = 0
$ws.Range("D418").Value = 0
$ws.Range("A419").Value = 44493
$ws.Range("B419").Value = 0
$ws.Range("C419").Value = 0
$ws.Range("D419").Value = 0
$ws.Range("A420").Value = 44494
$ws.Range("B420").Value = 0
$ws.Range("C420").Value = 0
$ws.Range("D420").Value = 0
$ws.Range("A421").Value = 44495
$ws.Range("B421").Value = 0
$ws.Range("C421").Value = 0
$ws.Range("D421").Value = 0
$ws.Range("A422").Value = 44496
$ws.Range("B422").Value = 0
$ws.Range("C422").Value = 0
$ws.Range("D422").Value = 0
$ws.Range("A423").Value = 44497
$ws.Range("B423").Value = 0
$ws.Range("C423").Value = 0
$ws.Range("D423").Value = 0
$ws.Range("A424").Value = 44498
$ws.Range("B424").Value = 0
$ws.Range("C424").Value = 0
$ws.Range("D424").Value = 0
$ws.Range("A425").Value = 44499
$ws.Range("B425").Value = 0
$ws.Range("C425").Value = 0
$ws.Range("D425").Value = 0
$ws.Range("A426").Value = 44500
$ws.Range("B426").Value = 0
$ws.Range("C426").Value = 0
$ws.Range("D426").Value = 0
$ws.Range("A427").Value = 44501
$ws.Range("B427").Value = 0
$ws.Range("C427").Value = 0
$ws.Range("D427").Value = 0
$ws.Range("A428").Value = 44502
$ws.Range("B428").Value = 0
$ws.Range("C428").Value = 0
$ws.Range("D428").Value = 0
$ws.Range("A429").Value = 44503
$ws.Range("B429").Value = 0
$ws.Range("C429").Value = 0
$ws.Range("D429").Value = 0
$ws.Range("A430").Value = 44504
$ws.Range("B430").Value = 0
$ws.Range("C430").Value = 0
$ws.Range("D430").Value = 0
$ws.Range("A431").Value = 44505
$ws.Range("B431").Value = 0
$ws.Range("C431").Value = 0
$ws.Range("D431").Value = 0
$ws.Range("A432").Value = 44506
$ws.Range("B432").Value = 0
$ws.Range("C432").Value = 0
$ws.Range("D432").Value = 0
$ws.Range("A433").Value = 44507
$ws.Range("B433").Value = 0
$ws.Range("C433").Value = 0
$ws.Range("D433").Value = 0
$ws.Range("A434").Value = 44508
$ws.Range("B434").Value = 0
$ws.Range("C434").Value = 0
$ws.Range("D434").Value = 0
$ws.Range("A435").Value = 44509
$ws.Range("B435").Value = 0
$ws.Range("C435").Value = 0
$ws.Range("D435").Value = 0
$ws.Range("A436").Value = 44510
$ws.Range("B436").Value = 0
$ws.Range("C436").Value = 0
$ws.Range("D436").Value = 0
$ws.Range("A437").Value = 44511
$ws.Range("B437").Value = 0
$ws.Range("C437").Value = 0
$ws.Range("D437").Value = 0
$ws.Range("A438").Value = 44512
$ws.Range("B438").Value = 0
$ws.Range("C438").Value = 0
$ws.Range("D438").Value = 0
$ws.Range("A439").Value = 44513
$ws.Range("B439").Value = 0
$ws.Range("C439").Value = 0
$ws.Range("D439").Value = 0
$ws.Range("A440").Value = 44514
$ws.Range("B440").Value = 0
$ws.Range("C440").Value = 0
$ws.Range("D440").Value = 0
$ws.Range("A441").Value = 44515
$ws.Range("B441").Value = 0
$ws.Range("C441").Value = 0
$ws.Range("D441").Value = 0
$ws.Range("A442").Value = 44516
$ws.Range("B442").Value = 1
$ws.Range("C442").Value = 1
$ws.Range("D442").Value = 43.78283712784589
$ws.Range("A443").Value = 44517
$ws.Range("B443").Value = 0
$ws.Range("C443").Value = 1
$ws.Range("D443").Value = 43.78283712784589
$ws.Range("A444").Value = 44518
$ws.Range("B444").Value = 0
$ws.Range("C444").Value = 1
$ws.Range("D444").Value = 43.78283712784589
$ws.Range("A445").Value = 44519
$ws.Range("B445").Value = 0
$ws.Range("C445").Value = 1
$ws.Range("D445").Value = 43.78283712784589
$ws.Range("A446").Value = 44520
$ws.Range("B446").Value = 0
$ws.Range("C446").Value = 1
$ws.Range("D446").Value = 43.78283712784589
$ws.Range("A447").Value = 44521
$ws.Range("B447").Value = 0
$ws.Range("C447").Value = 1
$ws.Range("D447").Value = 43.78283712784589
$ws.Range("A448").Value = 44522
$ws.Range("B448").Value = 0
$ws.Range("C448").Value = 1
$ws.Range("D448").Value = 43.78283712784589
$ws.Range("A449").Value = 44523
$ws.Range("B449").Value = 0
$ws.Range("C449").Value = 0
$ws.Range("D449").Value = 0
$ws.Range("A450").Value = 44524
$ws.Range("B450").Value = 1
$ws.Range("C450").Value = 1
$ws.Range("D450").Value = 43.78283712784589
$ws.Range("A451").Value = 44525
$ws.Range("B451").Value = 0
$ws.Range("C451").Value = 1
$ws.Range("D451").Value = 43.78283712784589
$ws.Range("A452").Value = 44526
$ws.Range("B452").Value = 0
$ws.Range("C452").Value = 1
$ws.Range("D452").Value = 43.78283712784589
$ws.Range("A453").Value = 44527
$ws.Range("B453").Value = 0
$ws.Range("C453").Value = 1
$ws.Range("D453").Value = 43.78283712784589
$ws.Range("A454").Value = 44528
$ws.Range("B454").Value = 0
$ws.Range("C454").Value = 1
$ws.Range("D454").Value = 43.78283712784589
$ws.Range("A455").Value = 44529
$ws.Range("B455").Value = 0
$ws.Range("C455").Value = 1
$ws.Range("D455").Value = 43.78283712784589
$ws.Range("A456").Value = 44530
$ws.Range("B456").Value = 0
$ws.Range("C456").Value = 1
$ws.Range("D456").Value = 43.78283712784589
$ws.Range("A457").Value = 44531
$ws.Range("B457").Value = 0
$ws.Range("C457").Value = 0
$ws.Range("D457").Value = 0
$ws.Range("A458").Value = 44532
$ws.Range("B458").Value = 1
$ws.Range("C458").Value = 1
$ws.Range("D458").Value = 43.78283712784589
$ws.Range("A459").Value = 44533
$ws.Range("B459").Value = 1
$ws.Range("C459").Value = 2
$ws.Range("D459").Value = 87.56567425569177
$ws.Range("A460").Value = 44534
$ws.Range("B460").Value = 0
$ws.Range("C460").Value = 2
$ws.Range("D460").Value = 87.56567425569177
$ws.Range("A461").Value = 44535
$ws.Range("B461").Value = 2
$ws.Range("C461").Value = 4
$ws.Range("D461").Value = 175.1313485113835
$ws.Range("A462").Value = 44536
$ws.Range("B462").Value = 0
$ws.Range("C462").Value = 4
$ws.Range("D462").Value = 175.1313485113835
$ws.Range("A463").Value = 44537
$ws.Range("B463").Value = 0
$ws.Range("C463").Value = 4
$ws.Range("D463").Value = 175.1313485113835
$ws.Range("A464").Value = 44538
$ws.Range("B464").Value = 1
$ws.Range("C464").Value = 5
$ws.Range("D464").Value = 218.9141856392294

Write-Output "Appended rows 386-464 (through 2021-12-08) to Sheet1."
